# Update column F ("dSF") values on Sheet1 to reflect repulled/recalculated data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = -8
$ws.Range("F3").Value = -5
$ws.Range("F4").Value = -4
$ws.Range("F5").Value = 4
$ws.Range("F6").Value = -3
$ws.Range("F7").Value = -5
$ws.Range("F8").Value = 11
